$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date (column C) for rows 2 through 13 from 2023-10-08 to 2023-10-09
# Use the raw Excel date serial number (45208) so no time-of-day fraction is introduced.
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 45208
}
